$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values following a recalculation/repull of data.
# Row -> new F value
$updates = @{
    6  = -3
    7  = -5
    8  = -2
    10 = -5
    12 = -2
    14 = -3
    15 = -1
    25 = -7
    26 = -3
    27 = 0
    33 = 2
    34 = -5
    40 = 0
    48 = 0
    49 = 0
    50 = 4
    53 = 5
    55 = -3
    57 = 3
    64 = 3
    66 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
